# Insert a new data row at row 54 (pushing existing rows 54-143 down to 55-144)
# and populate it with the new observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("54:54").Insert()

$ws.Range("A54").Value = 11
$ws.Range("B54").Value = "Vega Monumental Concepción"
$ws.Range("C54").Value = "Bíobío"
$ws.Range("D54").Value = 44757
$ws.Range("E54").Value = 8
$ws.Range("F54").Value = 100112043
$ws.Range("G54").Value = "Pepino ensalada"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 100
$ws.Range("K54").Value = 20000
$ws.Range("L54").Value = 22000
$ws.Range("M54").Value = 21000
$ws.Range("N54").Value = "`$/caja 60 unidades"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 350
$ws.Range("Q54").Value = 60
$ws.Range("R54").Value = "Hortaliza"
